$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Some Price values are plain numeric-looking strings (e.g. "505.55"); for those
# we force the cell to Text format first so they are stored as text (matching the
# source data which keeps prices as literal strings, not floating point numbers).

$ws.Range('D2').Value = '56.589.89'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '2.389.22'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '505.55'
$ws.Range('E5').Value = '  +1.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.38'
$ws.Range('E6').Value = '  +3.80%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.552'
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').Value = '2.393.81'
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0974'
$ws.Range('E10').Value = '  +2.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.150'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.334'
$ws.Range('E12').Value = '  +5.33%  '
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').Value = '2.812.89'
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('D15').Value = '56.555.15'
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.76'
$ws.Range('E16').Value = '  +1.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000133'
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('D18').Value = '2.369.51'
$ws.Range('E18').Value = '  -2.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.16'
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('E20').Value = '  +0.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '310.02'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.65'
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.40'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.370'
$ws.Range('E28').Value = '  -0.88%  '
$ws.Range('E29').Value = '  +1.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '174.05'
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('D31').Value = '0.0₃0726'
$ws.Range('E31').Value = '  +2.06%  '
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.12'
$ws.Range('E33').Value = '  +1.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.86'
$ws.Range('E34').Value = '  -3.80%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.88'
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('E38').Value = '  -0.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.84'
$ws.Range('E39').Value = '  +1.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.67'
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('E41').Value = '  +3.54%  '
$ws.Range('E42').Value = '  +0.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '131.62'
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.39'
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.83'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.567'
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0910'
$ws.Range('E47').Value = '  +1.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '246.93'
$ws.Range('E48').Value = '  -2.01%  '
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('E50').Value = '  +1.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.15'
$ws.Range('E51').Value = '  +6.49%  '
